$wb = $excel.ActiveWorkbook

# --- Repayment schedule sheet: insert a new (blank) column before column N ---
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Columns("N:N").Insert()

# Give the freshly inserted column N its own width (it currently just
# inherited column M's width/format from the insert).
$ws.Columns("N:N").ColumnWidth = 10.7109375

# Select the new active cell on this sheet and make it the active sheet/tab.
$ws.Range("R7").Select()
$ws.Activate()

# --- Acc_Repayment1 sheet: it is no longer the active/selected tab ---
$ws8 = $wb.Worksheets.Item("Acc_Repayment1")
$ws8.Range("E4").Select()
